$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# The document currently ends with a single "Section Headnote" paragraph
# that contains every resource's content concatenated into one big run of
# text. We split that into the proper sequence of resource/section
# paragraphs (with their own styles, numbering and TOC bookmarks), and
# shorten the original paragraph down to just the section's own headnote
# text ("What is a corporation?").
# -------------------------------------------------------------------------

# Step 1: locate the big trailing paragraph via a unique substring near its
# end, and replace its text (but keep its paragraph mark/style) with the
# new, short headnote text.
$rng = $d.Content
$rng.Find.Execute("2Section Two", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetPara = $rng.Paragraphs(1)
$targetRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)
$targetRange.Text = "What is a corporation?"

# Step 2: append 11 new (empty) paragraphs after it - one per new resource
# / section element we need to create.
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertAfter("`r`r`r`r`r`r`r`r`r`r`r")

$count = $d.Paragraphs.Count
$base = $count - 11

# 1.1 - Case of the District Number 1 -------------------------------------
$p = $d.Paragraphs($base + 1)
$p.Style = "Resource Number"
$p.Range.Text = "1.1"
$s = $p.Range.Start; $e = $p.Range.End
$d.Bookmarks.Add("_auto_toc_2", $d.Range($s, $e)) | Out-Null

$p = $d.Paragraphs($base + 2)
$p.Style = "Resource Title"
$p.Range.Text = "Case of the District Number 1"

$p = $d.Paragraphs($base + 3)
$p.Style = "Resource Headnote"

$p = $d.Paragraphs($base + 4)
$p.Style = "Case Text"
$p.Range.Text = "This is the body of case 1."

# 1.2 - Case of the District Number 2 -------------------------------------
$p = $d.Paragraphs($base + 5)
$p.Style = "Resource Number"
$p.Range.Text = "1.2"
$s = $p.Range.Start; $e = $p.Range.End
$d.Bookmarks.Add("_auto_toc_3", $d.Range($s, $e)) | Out-Null

$p = $d.Paragraphs($base + 6)
$p.Style = "Resource Title"
$p.Range.Text = "Case of the District Number 2"

$p = $d.Paragraphs($base + 7)
$p.Style = "Resource Headnote"
$p.Range.Text = "This is an annotatable resource in the casebook.`n"

$p = $d.Paragraphs($base + 8)
$p.Style = "Case Text"
$p.Range.Text = "highlighted: content to highlight; elided: content to elide; replaced: content to replace; commented: content to comment; highlighted2: second highlight content;`n"

# 2 - Section Two -----------------------------------------------------------
$p = $d.Paragraphs($base + 9)
$p.Style = "Section Number"
$p.Range.Text = "2"
$s = $p.Range.Start; $e = $p.Range.End
$d.Bookmarks.Add("_auto_toc_4", $d.Range($s, $e)) | Out-Null

$p = $d.Paragraphs($base + 10)
$p.Style = "Section Title"
$p.Range.Text = "Section Two"

$p = $d.Paragraphs($base + 11)
$p.Style = "Section Headnote"
$p.Range.Text = "This is the second chapter of the casebook.`n"
